$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 7509.5
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 7509.5
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 7509.5
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -8445.5
$ws.Range("H23").Value = 7509.5
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 7509.5
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 7509.5
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -7977.5
$ws.Range("H38").Value = 1140.7778
$ws.Range("I38").Value = 44.666668
$ws.Range("K38").Value = 134.000004
$ws.Range("M38").Value = 237.999996
$ws.Range("H42").Value = 98.875
$ws.Range("I42").Value = 98.875
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 296.625
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -66.625
$ws.Range("N42").ClearContents()
$ws.Range("H52").Value = 2659.6667
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 2659.6667
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 7979.000100000001
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -8299.000100000001
$ws.Range("H69").Value = 4331.6665
$ws.Range("J69").Value = 3995
$ws.Range("L69").Value = 11985
$ws.Range("N69").Value = -13733
$ws.Range("H70").Value = 4856.857
$ws.Range("J70").Value = 4856.857
$ws.Range("L70").Value = 14570.571
$ws.Range("N70").Value = -15110.571
$ws.Range("H72").Value = 4331.6665
$ws.Range("J72").Value = 3995
$ws.Range("L72").Value = 35955
$ws.Range("N72").Value = -44691
$ws.Range("H73").Value = 4856.857
$ws.Range("J73").Value = 4856.857
$ws.Range("L73").Value = 14570.571
$ws.Range("N73").Value = -16442.571
$ws.Range("H80").Value = 599.9286
$ws.Range("I80").Value = 599.9286
$ws.Range("K80").Value = 1799.7858
$ws.Range("M80").Value = -801.7857999999999
$ws.Range("H83").Value = 599.9286
$ws.Range("I83").Value = 599.9286
$ws.Range("K83").Value = 5399.3574
$ws.Range("M83").Value = -407.3573999999999
$ws.Range("H115").Value = 426
$ws.Range("I115").Value = 426
$ws.Range("K115").Value = 1278
$ws.Range("M115").Value = 289
$ws.Range("H116").Value = 6613.25
$ws.Range("I116").Value = 5770
$ws.Range("K116").Value = 5770
$ws.Range("M116").Value = -2328
$ws.Range("H118").Value = 606
$ws.Range("I118").Value = 553.7692
$ws.Range("K118").Value = 1661.3076
$ws.Range("M118").Value = -4.307599999999866
$ws.Range("H125").Value = 5623.125
$ws.Range("I125").Value = 2390
$ws.Range("K125").Value = 21510
$ws.Range("M125").Value = -19050
$ws.Range("H132").Value = 2570.5454
$ws.Range("I132").Value = 2570.5454
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7711.6362
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5181.6362
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 2305.1177
$ws.Range("J137").Value = 4549.25
$ws.Range("L137").Value = 13647.75
$ws.Range("N137").Value = -18747.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1626.1428
$ws.Range("I94").Value = 1849.3529
$ws.Range("K94").Value = 1849.3529
$ws.Range("M94").Value = -1398.3529
$ws.Range("H99").Value = 1149.6666
$ws.Range("I99").Value = 1157.6
$ws.Range("K99").Value = 1157.6
$ws.Range("M99").Value = 340.4000000000001
$ws.Range("H107").Value = 865.75
$ws.Range("I107").Value = 900.1667
$ws.Range("J107").Value = 762.5
$ws.Range("K107").Value = 900.1667
$ws.Range("L107").Value = 762.5
$ws.Range("M107").Value = 1019.8333
$ws.Range("N107").Value = -4602.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7108.2
$ws.Range("I58").Value = 4251.9165
$ws.Range("K58").Value = 4251.9165
$ws.Range("M58").Value = -4048.9165
$ws.Range("H99").Value = 7916.3335
$ws.Range("I99").Value = 7916.3335
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 7916.3335
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -6418.3335
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 997.6
$ws.Range("I107").Value = 1026
$ws.Range("K107").Value = 1026
$ws.Range("M107").Value = 894
$ws.Range("H126").Value = 7916.3335
$ws.Range("I126").Value = 7916.3335
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 23749.0005
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -21279.0005
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 7108.2
$ws.Range("I136").Value = 4251.9165
$ws.Range("K136").Value = 12755.7495
$ws.Range("M136").Value = -10205.7495

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 655
$ws.Range("J23").Value = 435
$ws.Range("L23").Value = 1305
$ws.Range("N23").Value = -1775
$ws.Range("H131").Value = 1600
$ws.Range("J131").Value = 1600
$ws.Range("L131").Value = 4800
$ws.Range("N131").Value = -14880

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 232.5
$ws.Range("I2").Value = 227.41667
$ws.Range("K2").Value = 227.41667
$ws.Range("M2").Value = -114.41667
$ws.Range("H92").Value = 7206.857
$ws.Range("J92").Value = 7206.857
$ws.Range("L92").Value = 7206.857
$ws.Range("N92").Value = -10950.857
$ws.Range("H132").Value = 2464.5
$ws.Range("I132").Value = 2129.6667
$ws.Range("J132").Value = 2799.3333
$ws.Range("K132").Value = 6389.000100000001
$ws.Range("L132").Value = 8397.999899999999
$ws.Range("M132").Value = -3859.000100000001
$ws.Range("N132").Value = -13457.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1975.5385
$ws.Range("I7").Value = 1516.5454
$ws.Range("J7").Value = 4500
$ws.Range("K7").Value = 1516.5454
$ws.Range("L7").Value = 4500
$ws.Range("M7").Value = -1404.5454
$ws.Range("N7").Value = -4724
$ws.Range("H22").Value = 1756.8182
$ws.Range("I22").Value = 1696.875
$ws.Range("J22").Value = 1916.6666
$ws.Range("K22").Value = 1696.875
$ws.Range("L22").Value = 1916.6666
$ws.Range("M22").Value = -1401.875
$ws.Range("N22").Value = -2506.6666
$ws.Range("H27").Value = 1756.8182
$ws.Range("I27").Value = 1696.875
$ws.Range("J27").Value = 1916.6666
$ws.Range("K27").Value = 1696.875
$ws.Range("L27").Value = 1916.6666
$ws.Range("M27").Value = -1589.875
$ws.Range("N27").Value = -2130.6666
$ws.Range("H40").Value = 7661.357
$ws.Range("I40").Value = 6362
$ws.Range("K40").Value = 6362
$ws.Range("M40").Value = -6226
$ws.Range("H126").Value = 1975.5385
$ws.Range("I126").Value = 1516.5454
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 4549.6362
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -2079.6362
$ws.Range("N126").Value = -18440
$ws.Range("H132").Value = 7933.1113
$ws.Range("I132").Value = 7342.5713
$ws.Range("K132").Value = 22027.7139
$ws.Range("M132").Value = -19497.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 751.38464
$ws.Range("I81").Value = 702
$ws.Range("J81").Value = 862.5
$ws.Range("K81").Value = 1404
$ws.Range("L81").Value = 1725
$ws.Range("M81").Value = -343
$ws.Range("N81").Value = -3847
$ws.Range("H84").Value = 751.38464
$ws.Range("I84").Value = 702
$ws.Range("J84").Value = 862.5
$ws.Range("K84").Value = 7020
$ws.Range("L84").Value = 8625
$ws.Range("M84").Value = -1716
$ws.Range("N84").Value = -19233
$ws.Range("H94").Value = 19145
$ws.Range("J94").Value = 19145
$ws.Range("L94").Value = 19145
$ws.Range("N94").Value = -20947
$ws.Range("H95").Value = 27034.5
$ws.Range("J95").Value = 27034.5
$ws.Range("L95").Value = 27034.5
$ws.Range("N95").Value = -32526.5
